$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.021.94"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.285.84"
$ws.Range("E3").Value = "  +2.67%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'495.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.54%  "
$ws.Range("D6").Value = "'127.82"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("D7").Value = "'0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "'0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.47%  "
$ws.Range("D9").Value = "2.283.95"
$ws.Range("E9").Value = "  +2.10%  "
$ws.Range("D10").Value = "'0.0949"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.77%  "
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("D12").Value = "'0.327"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.05%  "
$ws.Range("D13").Value = "'4.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.49%  "
$ws.Range("D14").Value = "2.664.09"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "'21.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.92%  "
$ws.Range("D16").Value = "53.998.99"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("D18").Value = "2.294.33"
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").Value = "'10.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +5.04%  "
$ws.Range("D20").Value = "'4.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.55%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'301.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.49%  "
$ws.Range("D23").Value = "'0.996"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("E24").Value = "  -2.67%  "
$ws.Range("D25").Value = "'62.49"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.24%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("E27").Value = "  +2.95%  "
$ws.Range("D28").Value = "2.380.81"
$ws.Range("E28").Value = "  +2.11%  "
$ws.Range("D29").Value = "'0.148"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.56%  "
$ws.Range("D31").Value = "'169.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "'1.61"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.36%  "
$ws.Range("E33").Value = "  +1.67%  "
$ws.Range("E34").Value = "  +2.77%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").Value = "'0.996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("D38").Value = "'17.67"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +4.47%  "
$ws.Range("D40").Value = "'0.866"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.29%  "
$ws.Range("D41").Value = "'3.70"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.98%  "
$ws.Range("D42").Value = "'35.42"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.22%  "
$ws.Range("E43").Value = "  +3.40%  "
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("D45").Value = "'3.36"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'4.94"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'127.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.81%  "
$ws.Range("D48").Value = "'0.0888"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.39%  "
$ws.Range("D49").Value = "'0.543"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.61%  "
$ws.Range("D50").Value = "'238.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.29%  "
$ws.Range("E51").Value = "  +3.03%  "
